# Trade #111 closed at 2026-02-17 16:02:42 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.92
$summary.Range("B4").Value = -1.09
$summary.Range("B5").Value = -0.2
$summary.Range("B6").Value = 111
$summary.Range("B8").Value = 55
$summary.Range("B9").Value = 36.04

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.92
$status.Range("D4").Value = 111
$status.Range("E4").Value = -1.09
$status.Range("F4").Value = -1.08
$status.Range("G4").Value = 36.04

# ---- Append the new closed trade (Trade #111) to the trade logs ----
function Add-TradeRow($ws) {
    $ws.Cells.Item(112, 1).Value = 111

    # Force the date-like string to stay as text instead of auto-converting
    # to a date serial number (matches how the rest of the log stores it).
    $dateCell = $ws.Cells.Item(112, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item(112, 3).Value = "16:02:36"
    $ws.Cells.Item(112, 4).Value = "MarketMaking"
    $ws.Cells.Item(112, 5).Value = "UP"
    $ws.Cells.Item(112, 6).Value = 0.87
    $ws.Cells.Item(112, 7).Value = 0.78
    $ws.Cells.Item(112, 8).Value = "CLOSED"
    $ws.Cells.Item(112, 9).Value = -10.3448
    $ws.Cells.Item(112, 10).Value = -0.09
    $ws.Cells.Item(112, 11).Value = 98.92
    $ws.Cells.Item(112, 12).Value = 0
    $ws.Cells.Item(112, 13).Value = 0
    $ws.Cells.Item(112, 14).Value = 0.6
    $ws.Cells.Item(112, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(112, 16).Value = "early_exit"
    $ws.Cells.Item(112, 17).Value = 0.14
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
